$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Insert a stand-alone "F" in the middle of the long run so the
#    run is split into three pieces (before / "F" / after), each
#    keeping its own (empty) run properties. Do this FIRST, while
#    absolute character offsets still match the pristine document.
# -----------------------------------------------------------------
$insertionPoint = 424
$ins = $d.Range($insertionPoint, $insertionPoint)
$ins.InsertAfter("F")

# Force the newly inserted character to become its own run by
# toggling a character property on it and back off again.
$newChar = $d.Range($insertionPoint, $insertionPoint + 1)
$newChar.Font.Bold = $true
$newChar.Font.Bold = $false

# Re-materialize explicit (empty) run properties on the two
# surrounding text runs so they match the original authoring style.
$beforeRun = $d.Range(66, $insertionPoint)
$beforeRun.Font.Bold = $true
$beforeRun.Font.Bold = $false

$afterRun = $d.Range($insertionPoint + 1, 833)
$afterRun.Font.Bold = $true
$afterRun.Font.Bold = $false

# -----------------------------------------------------------------
# 2) Simple single/unique text replacements (each target substring
#    occurs exactly once in the whole document, so Find/Replace is
#    safe and is insensitive to the offset shift caused by step 1).
# -----------------------------------------------------------------
$d.Content.Find.Execute("2", $true, $false, $false, $false, $false, $true, 1, $false, "FF", 2)
$d.Content.Find.Execute("1", $true, $false, $false, $false, $false, $true, 1, $false, "KNDK", 2)
$d.Content.Find.Execute("3", $true, $false, $false, $false, $false, $true, 1, $false, "FNDJNF", 2)
$d.Content.Find.Execute("5", $true, $false, $false, $false, $false, $true, 1, $false, "F", 2)

# -----------------------------------------------------------------
# 3) Normal style / document default paragraph formatting tweaks.
# -----------------------------------------------------------------
$normal = $d.Styles.Item("Normal")
$normal.ParagraphFormat.SpaceBefore = 0
$normal.ParagraphFormat.SpaceAfter = 0
$normal.ParagraphFormat.Alignment = 0
